$d = $word.ActiveDocument

# 1) Remove the existing "_GoBack" bookmark (it currently sits in the first
#    paragraph). It will be re-created further down, at the end of the
#    document, as part of the newly appended paragraphs.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2) Append four new paragraphs after the current last paragraph (the one
#    holding the last inline picture), just before the final section break.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$tail = $lastPara.Range
$tail.Collapse(0)
$tail.InsertParagraphAfter() | Out-Null

$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertionRange = $newPara.Range

$w = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$xml = ""
$xml += "<w:p $w>"
$xml +=   "<w:pPr>"
$xml +=     "<w:ind w:left='360'/>"
$xml +=     "<w:jc w:val='both'/>"
$xml +=     "<w:rPr>"
$xml +=       "<w:rFonts w:ascii='Verdana' w:hAnsi='Verdana'/>"
$xml +=       "<w:sz w:val='24'/>"
$xml +=       "<w:szCs w:val='24'/>"
$xml +=     "</w:rPr>"
$xml +=   "</w:pPr>"
$xml += "</w:p>"

$xml += "<w:p $w>"
$xml +=   "<w:pPr>"
$xml +=     "<w:ind w:left='360'/>"
$xml +=     "<w:jc w:val='center'/>"
$xml +=     "<w:rPr>"
$xml +=       "<w:rFonts w:ascii='Verdana' w:hAnsi='Verdana'/>"
$xml +=       "<w:b/>"
$xml +=       "<w:sz w:val='24'/>"
$xml +=       "<w:szCs w:val='24'/>"
$xml +=     "</w:rPr>"
$xml +=   "</w:pPr>"
$xml +=   "<w:r>"
$xml +=     "<w:rPr>"
$xml +=       "<w:rFonts w:ascii='Verdana' w:hAnsi='Verdana'/>"
$xml +=       "<w:b/>"
$xml +=       "<w:sz w:val='24'/>"
$xml +=       "<w:szCs w:val='24'/>"
$xml +=     "</w:rPr>"
$xml +=     "<w:t>Criando o serviço do jogo da velha</w:t>"
$xml +=   "</w:r>"
$xml += "</w:p>"

$xml += "<w:p $w>"
$xml +=   "<w:pPr>"
$xml +=     "<w:ind w:left='360'/>"
$xml +=     "<w:jc w:val='both'/>"
$xml +=     "<w:rPr>"
$xml +=       "<w:rFonts w:ascii='Verdana' w:hAnsi='Verdana'/>"
$xml +=       "<w:sz w:val='24'/>"
$xml +=       "<w:szCs w:val='24'/>"
$xml +=     "</w:rPr>"
$xml +=   "</w:pPr>"
$xml +=   "<w:r>"
$xml +=     "<w:rPr>"
$xml +=       "<w:rFonts w:ascii='Verdana' w:hAnsi='Verdana'/>"
$xml +=       "<w:sz w:val='24'/>"
$xml +=       "<w:szCs w:val='24'/>"
$xml +=     "</w:rPr>"
$xml +=     "<w:t>Criando o serviço do jogo da velha que conterá regras de negócio nosso jogo da velha vai ter regras para registrar uma vitória, uma jogada e a jogada computador será codificado nesse arquivo.</w:t>"
$xml +=   "</w:r>"
$xml += "</w:p>"

$xml += "<w:p $w>"
$xml +=   "<w:pPr>"
$xml +=     "<w:ind w:left='360'/>"
$xml +=     "<w:jc w:val='both'/>"
$xml +=     "<w:rPr>"
$xml +=       "<w:rFonts w:ascii='Verdana' w:hAnsi='Verdana'/>"
$xml +=       "<w:sz w:val='24'/>"
$xml +=       "<w:szCs w:val='24'/>"
$xml +=     "</w:rPr>"
$xml +=   "</w:pPr>"
$xml +=   "<w:bookmarkStart w:id='0' w:name='_GoBack'/>"
$xml +=   "<w:bookmarkEnd w:id='0'/>"
$xml += "</w:p>"

$insertionRange.InsertXML($xml) | Out-Null

Write-Host "Added service-section paragraphs and relocated the _GoBack bookmark."
